$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.544.33'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '2.601.68'
$ws.Range('E3').Value = '  +1.69%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '515.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.598'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.78%  '
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('E10').Value = '  +1.82%  '
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '3.056.85'
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('D14').Value = '60.559.88'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.63'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000141'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').Value = '2.602.44'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '358.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.51%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.428'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.70%  '
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').Value = '2.715.97'
$ws.Range('E26').Value = '  +0.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = '0.0₃0838'
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('E32').Value = '  +2.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.03'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.940'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +11.44%  '
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('E38').Value = '  +1.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.34'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.77%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.842'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '288.79'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('E43').Value = '  +2.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.619'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.98%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0558'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.997'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0236'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.97%  '
$ws.Range('E50').Value = '  +0.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.90%  '
